$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TTD")
Write-Output $ws.Name
